$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.794.46"
$ws.Range("E2").Value = "  -1.65%  "

$ws.Range("D3").Value = "2.562.94"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "302.38"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").Value = "95.22"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "0.545"

$ws.Range("D10").Value = "36.02"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("E13").Value = "  +6.36%  "

$ws.Range("D14").Value = "2.582.30"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").Value = "0.879"
$ws.Range("E15").Value = "  -1.48%  "

$ws.Range("E16").Value = "  -1.02%  "

$ws.Range("D17").Value = "42.832.74"
$ws.Range("E17").Value = "  -1.66%  "

$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("D21").Value = "71.80"
$ws.Range("E21").Value = "  -1.61%  "

$ws.Range("D22").Value = "253.15"
$ws.Range("E22").Value = "  -5.08%  "

$ws.Range("D23").Value = "2.94"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  -3.93%  "

$ws.Range("D25").Value = "28.67"

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "36.98"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "6.01"
$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("D31").Value = "154.36"
$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -4.02%  "

$ws.Range("D33").Value = "3.38"
$ws.Range("E33").Value = "  -6.79%  "

$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "0.0799"
$ws.Range("E35").Value = "  -1.99%  "

$ws.Range("D36").Value = "18.43"
$ws.Range("E36").Value = "  +10.01%  "

$ws.Range("E37").Value = "  -4.28%  "

$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").Value = "23.45"
$ws.Range("E39").Value = "  -4.54%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0311"
$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  -5.35%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "2.07"
$ws.Range("E42").Value = "  +28.66%  "

$ws.Range("D43").Value = "3.88"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").Value = "2.082.95"
$ws.Range("E44").Value = "  +1.97%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "9.24"
$ws.Range("E46").Value = "  +1.91%  "

$ws.Range("D47").Value = "85.01"
$ws.Range("E47").Value = "  -3.99%  "

$ws.Range("D48").Value = "75.94"
$ws.Range("E48").Value = "  +9.59%  "

$ws.Range("D49").Value = "106.45"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").Value = "2.817.03"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  +0.83%  "
